$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.1857778469365516
$ws.Range("D2").Value = 0.8537229047962624

$ws.Range("C3").Value = 0.01070189892367874
$ws.Range("D3").Value = 0.991523833272711

$ws.Range("C4").Value = 0.2805866942285271
$ws.Range("D4").Value = 0.7807277463332984

$ws.Range("C5").Value = -0.3306339371931933
$ws.Range("D5").Value = 0.7429499516626406

$ws.Range("C6").Value = 0.2599959592849306
$ws.Range("D6").Value = 0.796434935887649

$ws.Range("C7").Value = 0.9931572636137704
$ws.Range("D7").Value = 0.3276492542044294

$ws.Range("C8").Value = -0.2584335074202835
$ws.Range("D8").Value = 0.7976304338836364

$ws.Range("C9").Value = 0.3705418230099904
$ws.Range("D9").Value = 0.7132774029178419

$ws.Range("C10").Value = -0.4464825655076207
$ws.Range("D10").Value = 0.6580795915645936

$ws.Range("C11").Value = -0.6385527812838422
$ws.Range("D11").Value = 0.5273921692542036
